$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.916621666666667
$ws.Range("H2").Value = 17.749865
$ws.Range("I2").Value = 0.2364294176430088
$ws.Range("J2").Value = 0.2364294176430088
$ws.Range("M2").Value = 8.533046666666666
$ws.Range("N2").Value = 25.59914
$ws.Range("O2").Value = 0.2932132236642383
$ws.Range("P2").Value = 0.2932132236642383
$ws.Range("Q2").Value = 50.48680879067777
$ws.Range("R2").Value = 454.3812791161
$ws.Range("S2").Value = 0.06932423171616514
$ws.Range("T2").Value = 0.06932423171616514

# Row 3
$ws.Range("G3").Value = 5.916621666666667
$ws.Range("H3").Value = 17.749865
$ws.Range("I3").Value = 0.2364294176430088
$ws.Range("J3").Value = 0.2364294176430088
$ws.Range("O3").Value = 0.3119288965200195
$ws.Range("P3").Value = 0.3119288965200194
$ws.Range("Q3").Value = 53.7093598920589
$ws.Range("R3").Value = 483.38423902853
$ws.Range("S3").Value = 0.07374916735025455
$ws.Range("T3").Value = 0.07374916735025454

# Row 4
$ws.Range("G4").Value = 5.916621666666667
$ws.Range("H4").Value = 17.749865
$ws.Range("I4").Value = 0.2364294176430088
$ws.Range("J4").Value = 0.2364294176430088
$ws.Range("O4").Value = 0.3948578798157423
$ws.Range("P4").Value = 0.3948578798157423
$ws.Range("Q4").Value = 67.98845573410333
$ws.Range("R4").Value = 611.89610160693
$ws.Range("S4").Value = 0.09335601857658909
$ws.Range("T4").Value = 0.09335601857658909

# Row 5
$ws.Range("G5").Value = 9.915995333333335
$ws.Range("I5").Value = 0.3962452112189236
$ws.Range("J5").Value = 0.3962452112189236
$ws.Range("M5").Value = 8.533046666666666
$ws.Range("N5").Value = 25.59914
$ws.Range("O5").Value = 0.2932132236642383
$ws.Range("P5").Value = 0.2932132236642383
$ws.Range("Q5").Value = 84.61365092578222
$ws.Range("R5").Value = 761.5228583320401
$ws.Range("S5").Value = 0.1161843357430176
$ws.Range("T5").Value = 0.1161843357430176

# Row 6
$ws.Range("G6").Value = 9.915995333333335
$ws.Range("I6").Value = 0.3962452112189236
$ws.Range("J6").Value = 0.3962452112189236
$ws.Range("O6").Value = 0.3119288965200195
$ws.Range("P6").Value = 0.3119288965200194
$ws.Range("Q6").Value = 90.01450355469913
$ws.Range("R6").Value = 810.1305319922922
$ws.Range("S6").Value = 0.1236003314868609
$ws.Range("T6").Value = 0.1236003314868609

# Row 7
$ws.Range("G7").Value = 9.915995333333335
$ws.Range("I7").Value = 0.3962452112189236
$ws.Range("J7").Value = 0.3962452112189236
$ws.Range("O7").Value = 0.3948578798157423
$ws.Range("P7").Value = 0.3948578798157423
$ws.Range("Q7").Value = 113.9456344788947
$ws.Range("S7").Value = 0.1564605439890451
$ws.Range("T7").Value = 0.1564605439890452

# Row 8
$ws.Range("I8").Value = 0.3673253711380675
$ws.Range("J8").Value = 0.3673253711380676
$ws.Range("M8").Value = 8.533046666666666
$ws.Range("N8").Value = 25.59914
$ws.Range("O8").Value = 0.2932132236642383
$ws.Range("P8").Value = 0.2932132236642383
$ws.Range("Q8").Value = 78.43814852436888
$ws.Range("R8").Value = 705.94333671932
$ws.Range("S8").Value = 0.1077046562050555
$ws.Range("T8").Value = 0.1077046562050556

# Row 9
$ws.Range("I9").Value = 0.3673253711380675
$ws.Range("J9").Value = 0.3673253711380676
$ws.Range("O9").Value = 0.3119288965200195
$ws.Range("P9").Value = 0.3119288965200194
$ws.Range("R9").Value = 751.0033936282362
$ws.Range("S9").Value = 0.114579397682904
$ws.Range("T9").Value = 0.114579397682904

# Row 10
$ws.Range("I10").Value = 0.3673253711380675
$ws.Range("J10").Value = 0.3673253711380676
$ws.Range("O10").Value = 0.3948578798157423
$ws.Range("P10").Value = 0.3948578798157423
$ws.Range("R10").Value = 950.664113042316
$ws.Range("S10").Value = 0.145041317250108
$ws.Range("T10").Value = 0.145041317250108
